# Weekly data refresh: Fruta / hortaliza, semanal
# Updates Fecha, Volumen, Precio minimo/maximo/promedio ponderado,
# Origen and Precio $/Kg for rows 2-29 (Hortaliza - Haba @ Vega Monumental Concepcion).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 7 changes
$ws.Range("D2").Value = 44476
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 7500
$ws.Range("M2").Value = 7250
$ws.Range("O2").Value = 'Región Metropolitana'
$ws.Range("P2").Value = 290
# Row 3: 7 changes
$ws.Range("D3").Value = 44540
$ws.Range("J3").Value = 140
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11429
$ws.Range("O3").Value = 'Región del Maule'
$ws.Range("P3").Value = 457
# Row 4: 2 changes
$ws.Range("D4").Value = 44335
$ws.Range("O4").Value = 'Provincia de Limarí'
# Row 5: 7 changes
$ws.Range("D5").Value = 44537
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 8500
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 8719
$ws.Range("O5").Value = 'Región del Maule'
$ws.Range("P5").Value = 349
# Row 6: 7 changes
$ws.Range("D6").Value = 44545
$ws.Range("J6").Value = 140
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14429
$ws.Range("O6").Value = 'Provincia de Chacabuco'
$ws.Range("P6").Value = 577
# Row 7: 7 changes
$ws.Range("D7").Value = 44461
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 13500
$ws.Range("O7").Value = 'Provincia del Elquí'
$ws.Range("P7").Value = 540
# Row 8: 5 changes
$ws.Range("D8").Value = 44316
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17000
$ws.Range("P8").Value = 680
# Row 9: 6 changes
$ws.Range("D9").Value = 44526
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 7500
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 7750
$ws.Range("P9").Value = 310
# Row 10: 7 changes
$ws.Range("D10").Value = 44482
$ws.Range("J10").Value = 430
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 8500
$ws.Range("M10").Value = 8267
$ws.Range("O10").Value = 'Región de O''Higgins'
$ws.Range("P10").Value = 331
# Row 11: 5 changes
$ws.Range("D11").Value = 44188
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 19000
$ws.Range("P11").Value = 760
# Row 12: 6 changes
$ws.Range("D12").Value = 44523
$ws.Range("K12").Value = 9000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 9500
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 380
# Row 13: 7 changes
$ws.Range("D13").Value = 44755
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 16000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 16500
$ws.Range("O13").Value = 'Región de Coquimbo'
$ws.Range("P13").Value = 660
# Row 14: 1 changes
$ws.Range("D14").Value = 44160
# Row 15: 5 changes
$ws.Range("D15").Value = 44467
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 8500
$ws.Range("P15").Value = 340
# Row 16: 5 changes
$ws.Range("D16").Value = 44505
$ws.Range("J16").Value = 180
$ws.Range("M16").Value = 6222
$ws.Range("O16").Value = 'Región del Maule'
$ws.Range("P16").Value = 249
# Row 17: 6 changes
$ws.Range("D17").Value = 44498
$ws.Range("J17").Value = 220
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 7500
$ws.Range("M17").Value = 7273
$ws.Range("P17").Value = 291
# Row 18: 5 changes
$ws.Range("D18").Value = 44692
$ws.Range("K18").Value = 20000
$ws.Range("L18").Value = 22000
$ws.Range("M18").Value = 21000
$ws.Range("P18").Value = 840
# Row 19: 6 changes
$ws.Range("D19").Value = 44454
$ws.Range("K19").Value = 13000
$ws.Range("L19").Value = 14000
$ws.Range("M19").Value = 13500
$ws.Range("O19").Value = 'Provincia del Elquí'
$ws.Range("P19").Value = 540
# Row 20: 6 changes
$ws.Range("D20").Value = 44351
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = 15500
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("P20").Value = 620
# Row 21: 6 changes
$ws.Range("D21").Value = 44509
$ws.Range("K21").Value = 6500
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 6750
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 270
# Row 22: 7 changes
$ws.Range("D22").Value = 44533
$ws.Range("J22").Value = 180
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 8500
$ws.Range("M22").Value = 8222
$ws.Range("O22").Value = 'Región del Maule'
$ws.Range("P22").Value = 329
# Row 23: 6 changes
$ws.Range("D23").Value = 44673
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 18000
$ws.Range("L23").Value = 19000
$ws.Range("M23").Value = 18375
$ws.Range("P23").Value = 735
# Row 24: 7 changes
$ws.Range("D24").Value = 44503
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 9000
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 9400
$ws.Range("O24").Value = 'Provincia de Melipilla'
$ws.Range("P24").Value = 376
# Row 25: 7 changes
$ws.Range("D25").Value = 44384
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 13000
$ws.Range("M25").Value = 12500
$ws.Range("O25").Value = 'Región de Coquimbo'
$ws.Range("P25").Value = 500
# Row 26: 7 changes
$ws.Range("D26").Value = 44162
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 7500
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 7750
$ws.Range("O26").Value = 'Región Metropolitana'
$ws.Range("P26").Value = 310
# Row 27: 6 changes
$ws.Range("D27").Value = 44483
$ws.Range("J27").Value = 350
$ws.Range("K27").Value = 5500
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = 5714
$ws.Range("P27").Value = 229
# Row 28: 5 changes
$ws.Range("D28").Value = 44497
$ws.Range("J28").Value = 150
$ws.Range("M28").Value = 6333
$ws.Range("O28").Value = 'Región Metropolitana'
$ws.Range("P28").Value = 253
# Row 29: 7 changes
$ws.Range("D29").Value = 44517
$ws.Range("J29").Value = 130
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 6500
$ws.Range("M29").Value = 6269
$ws.Range("O29").Value = 'Región Metropolitana'
$ws.Range("P29").Value = 251